{"js": "// Add \"make all loaders one class\" as the text of the empty paragraph\n// that immediately follows the \"skips\" paragraph (the document's 2nd\n// paragraph overall), without inserting any new paragraphs.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items[1];\ntargetParagraph.insertText(\"make all loaders one class\", Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Add \"make all loaders one class\" as the text of the empty paragraph\n# that immediately follows the \"skips\" paragraph (the document's 2nd\n# paragraph overall), without inserting any new paragraphs.\n$d = $word.ActiveDocument\n$p = $d.Paragraphs.Item(2)\n$p.Range.InsertBefore(\"make all loaders one class\")\n"}
